# Bump the "Published" date on the cover page from
# "October 8, 2020" to "October 9, 2020" (the document's day-of-month
# typo/update — the only textual change in this revision).
#
# In the underlying OOXML this single digit is its own isolated run
# (preceded by a run ending in "October " and followed by a run
# starting with ", 2020"), so scoping the Find/Replace to the
# "Published: October 8, 2020" paragraph and matching the whole word
# "8" hits exactly that run and nothing else in the document.

$d = $word.ActiveDocument

$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Published:*") {
        $targetParagraph = $p
        break
    }
}

if ($targetParagraph -eq $null) {
    # Fallback: search the whole document if the expected paragraph
    # text wasn't found for some reason.
    $searchRange = $d.Content
} else {
    $searchRange = $targetParagraph.Range
}

$find = $searchRange.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$replaced = $find.Execute(
    "8",    # FindText
    $true,  # MatchWholeWord
    $false, # MatchCase
    $false, # MatchWildcards
    $false, # MatchSoundsLike
    $false, # MatchAllWordForms
    $true,  # Forward
    1,      # Wrap (wdFindContinue)
    $false, # Format
    "9",    # ReplaceWith
    2       # Replace (wdReplaceAll)
)

Write-Output ("Replaced: " + $replaced)
if ($targetParagraph -ne $null) {
    Write-Output ("Paragraph now reads: " + $targetParagraph.Range.Text)
}
